$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose numeric text values are changing. All of these are stored in
# the workbook as text (shared strings), not as numeric cells, so we force
# text number-formatting before writing the value (otherwise Excel's COM
# layer auto-detects the numeric-looking string and stores it as a real
# number), then clear the formatting again so the cell keeps the same
# "General"-looking appearance it started with.
$updates = @{
    "C2" = "-936722091.828046"
    "C3" = "850894453.187894"
    "C4" = "1744702725.69586"
    "C5" = "2280987689.20065"
    "D2" = "3732.83737464552"
    "D3" = "3732.83737464552"
    "D4" = "3732.83737464552"
    "D5" = "3732.83737464552"
    "E2" = "9586774.16410338"
    "E3" = "9586774.16410338"
    "E4" = "9586774.16410338"
    "E5" = "9586774.16410338"
    "F2" = "29500633.6890926"
    "F3" = "29500633.6890926"
    "F4" = "29500633.6890926"
    "F5" = "29500633.6890926"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).ClearFormats()
}
